# Reorder the "estado de cuenta" data rows so records are grouped by
# period (1808 first, then 1809) instead of by worker.
#
# Before:
#   16: 32721014 | ESTHER MARINA FRAGOSO LUBO       | 1809
#   17: 32721014 | ESTHER MARINA FRAGOSO LUBO       | 1808
#   18: 91519089 | HAILE YASSER CASTAÑEDA CAZES     | 1809
#   19: 91519089 | HAILE YASSER CASTAÑEDA CAZES     | 1808
#   20: 45761960 | CLAUDIA ISABEL RIPOLL BENAVIDES  | 1809
#   21: 45761960 | CLAUDIA ISABEL RIPOLL BENAVIDES  | 1808
#
# After:
#   16: 32721014 | ESTHER MARINA FRAGOSO LUBO       | 1808
#   17: 91519089 | HAILE YASSER CASTAÑEDA CAZES     | 1808
#   18: 45761960 | CLAUDIA ISABEL RIPOLL BENAVIDES  | 1808
#   19: 32721014 | ESTHER MARINA FRAGOSO LUBO       | 1809
#   20: 91519089 | HAILE YASSER CASTAÑEDA CAZES     | 1809
#   21: 45761960 | CLAUDIA ISABEL RIPOLL BENAVIDES  | 1809

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$docNumbers = @("32721014", "91519089", "45761960")
$names = @("ESTHER MARINA FRAGOSO LUBO", "HAILE YASSER CASTAÑEDA CAZES", "CLAUDIA ISABEL RIPOLL BENAVIDES")
$periods = @("1808", "1809")

$row = 16
foreach ($period in $periods) {
    for ($i = 0; $i -lt 3; $i++) {
        $ws.Cells.Item($row, 3).Value = $docNumbers[$i]
        $ws.Cells.Item($row, 4).Value = $names[$i]
        $ws.Cells.Item($row, 5).Value = $period
        $row++
    }
}
